{"js": "const replacements = [\n  [\"157\u00d75=785\", \"784\u00d77=5488\"],\n  [\"630\u00d75=3150\", \"892\u00d75=4460\"],\n  [\"333\u00d77=2331\", \"752\u00d78=6016\"],\n  [\"326\u00d73=978\", \"431\u00d79=3879\"],\n  [\"392\u00d76=2352\", \"489\u00d72=978\"],\n  [\"375\u00d78=3000\", \"610\u00d76=3660\"],\n  [\"334\u00d76=2004\", \"257\u00d74=1028\"],\n  [\"174\u00d78=1392\", \"445\u00d76=2670\"],\n  [\"739\u00d77=5173\", \"427\u00d72=854\"],\n  [\"932\u00d79=8388\", \"131\u00d79=1179\"],\n  [\"988\u00d76=5928\", \"984\u00d74=3936\"],\n  [\"117\u00d76=702\", \"209\u00d72=418\"],\n  [\"983\u00d78=7864\", \"714\u00d72=1428\"],\n  [\"914\u00d79=8226\", \"665\u00d76=3990\"],\n  [\"544\u00d73=1632\", \"700\u00d72=1400\"],\n  [\"204\u00d76=1224\", \"175\u00d72=350\"],\n  [\"495\u00d73=1485\", \"911\u00d73=2733\"],\n  [\"102\u00d72=204\", \"650\u00d75=3250\"],\n  [\"901\u00d79=8109\", \"359\u00d78=2872\"],\n  [\"815\u00d73=2445\", \"805\u00d78=6440\"],\n  [\"997\u00d76=5982\", \"923\u00d76=5538\"],\n  [\"238\u00d73=714\", \"796\u00d77=5572\"],\n  [\"214\u00d77=1498\", \"722\u00d79=6498\"],\n  [\"874\u00d79=7866\", \"513\u00d79=4617\"],\n  [\"678\u00d75=3390\", \"902\u00d75=4510\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @('157\u00d75=785', '784\u00d77=5488'),\n  @('630\u00d75=3150', '892\u00d75=4460'),\n  @('333\u00d77=2331', '752\u00d78=6016'),\n  @('326\u00d73=978', '431\u00d79=3879'),\n  @('392\u00d76=2352', '489\u00d72=978'),\n  @('375\u00d78=3000', '610\u00d76=3660'),\n  @('334\u00d76=2004', '257\u00d74=1028'),\n  @('174\u00d78=1392', '445\u00d76=2670'),\n  @('739\u00d77=5173', '427\u00d72=854'),\n  @('932\u00d79=8388', '131\u00d79=1179'),\n  @('988\u00d76=5928', '984\u00d74=3936'),\n  @('117\u00d76=702', '209\u00d72=418'),\n  @('983\u00d78=7864', '714\u00d72=1428'),\n  @('914\u00d79=8226', '665\u00d76=3990'),\n  @('544\u00d73=1632', '700\u00d72=1400'),\n  @('204\u00d76=1224', '175\u00d72=350'),\n  @('495\u00d73=1485', '911\u00d73=2733'),\n  @('102\u00d72=204', '650\u00d75=3250'),\n  @('901\u00d79=8109', '359\u00d78=2872'),\n  @('815\u00d73=2445', '805\u00d78=6440'),\n  @('997\u00d76=5982', '923\u00d76=5538'),\n  @('238\u00d73=714', '796\u00d77=5572'),\n  @('214\u00d77=1498', '722\u00d79=6498'),\n  @('874\u00d79=7866', '513\u00d79=4617'),\n  @('678\u00d75=3390', '902\u00d75=4510'),\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $found = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $found) {\n        throw \"No match found for: $old\"\n    }\n}\n"}
